# "Generate Report for Handback"
#
# This script updates the localization-status workbook to reflect that the
# de-de handback has completed (and zh-cn's handback report/target file
# columns get populated), matching the "Generate Report for Handback" run.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------------
# Overview sheet: status text for both rows changes from "Ready for handoff"
# to "Handed back: in sync with en-US" (columns E = zh-cn, F = de-de).
# ---------------------------------------------------------------------------
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"

# Widen the status columns so the longer text fits.
$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668

# ---------------------------------------------------------------------------
# zh-cn sheet: populate "Latest Target File" (I) and "Latest Handback File"
# (J) for both tracked files.
# ---------------------------------------------------------------------------
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("I2"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5dc06455c21191739bf3602547d1f13809ebac6a/e2e/b379c699-84d7-46ea-9aa5-29e98526f75f.md",
    "",
    "",
    "b379c699-84d7-46ea-9aa5-29e98526f75f.md"
) | Out-Null
$wsZhCn.Range("J2").Value = "b379c699-84d7-46ea-9aa5-29e98526f75f.5a10a56527eb346c3e54a9c8a6a25ef99fde7fd5.zh-cn.xlf"

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("I3"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5dc06455c21191739bf3602547d1f13809ebac6a/e2e/b6e3eb6c-30d7-44cb-ac95-914adee7e921.md",
    "",
    "",
    "b6e3eb6c-30d7-44cb-ac95-914adee7e921.md"
) | Out-Null
$wsZhCn.Range("J3").Value = "b6e3eb6c-30d7-44cb-ac95-914adee7e921.66b8fb4a60aba6712040b3c75130476f2e1ea5b2.zh-cn.xlf"


# "Latest Handback DateTime" (K) was an empty 0001-01-01 placeholder and is
# now filled in with the real handback timestamp for both rows.
$wsZhCn.Range("K2").Value = "2016-09-01 20:55:42"
$wsZhCn.Range("K3").Value = "2016-09-01 20:55:42"

$wsZhCn.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsZhCn.Columns.Item(9).ColumnWidth = 39.166666666666664
$wsZhCn.Columns.Item(10).ColumnWidth = 39.166666666666664

# ---------------------------------------------------------------------------
# de-de sheet: same "Latest Target File" / "Latest Handback File" columns,
# plus the "Latest Handback DateTime" (K) now has a real handback time
# instead of the 0001-01-01 00:00:00 placeholder.
# ---------------------------------------------------------------------------
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("I2"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5dc06455c21191739bf3602547d1f13809ebac6a/e2e/b379c699-84d7-46ea-9aa5-29e98526f75f.md",
    "",
    "",
    "b379c699-84d7-46ea-9aa5-29e98526f75f.md"
) | Out-Null
$wsDeDe.Range("J2").Value = "b379c699-84d7-46ea-9aa5-29e98526f75f.5a10a56527eb346c3e54a9c8a6a25ef99fde7fd5.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-09-01 20:55:50"

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("I3"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5dc06455c21191739bf3602547d1f13809ebac6a/e2e/b6e3eb6c-30d7-44cb-ac95-914adee7e921.md",
    "",
    "",
    "b6e3eb6c-30d7-44cb-ac95-914adee7e921.md"
) | Out-Null
$wsDeDe.Range("J3").Value = "b6e3eb6c-30d7-44cb-ac95-914adee7e921.66b8fb4a60aba6712040b3c75130476f2e1ea5b2.de-de.xlf"
$wsDeDe.Range("K3").Value = "2016-09-01 20:55:50"

$wsDeDe.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsDeDe.Columns.Item(9).ColumnWidth = 39.166666666666664
$wsDeDe.Columns.Item(10).ColumnWidth = 39.166666666666664
